$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume change (E) columns with latest scraped values.
# D-column values are textual price strings (may look numeric) so force text format
# to avoid Excel auto-converting them to numbers and losing formatting (e.g. trailing zeros).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.085.83'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.651.82'
$ws.Range('E3').Value = '  -0.73%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.06'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5274'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -1.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06311'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('E10').Value = '  -2.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07795'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.518'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.665.76'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.878.39'
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5469'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8187'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.30'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.078.73'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.576'
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.47'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.05'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.68'
$ws.Range('E25').Value = '  +3.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1230'
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.209'
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.451'
$ws.Range('E29').Value = '  +2.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05801'
$ws.Range('E30').Value = '  -2.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.273'
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.258'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.413'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  -1.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5743'
$ws.Range('E38').Value = '  +1.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01609'
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8497'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '104.33'
$ws.Range('E41').Value = '  +3.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.711'
$ws.Range('E43').Value = '  -4.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.029.45'
$ws.Range('E44').Value = '  +2.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.792.90'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.95'
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4329'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.840'
$ws.Range('E49').Value = '  -2.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05142'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.445'
$ws.Range('E51').Value = '  -1.40%  '
